# Rotate columns D,E,F,G right by one position on every row of the sheet.
# i.e. for each row: newD = oldG, newE = oldD, newF = oldE, newG = oldF
# (A,B,C are left untouched.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $d = $ws.Cells.Item($r, 4).Value()
    $e = $ws.Cells.Item($r, 5).Value()
    $f = $ws.Cells.Item($r, 6).Value()
    $g = $ws.Cells.Item($r, 7).Value()

    $ws.Cells.Item($r, 4).Value = $g
    $ws.Cells.Item($r, 5).Value = $d
    $ws.Cells.Item($r, 6).Value = $e
    $ws.Cells.Item($r, 7).Value = $f
}
